$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as TEXT (shared string) even when it looks
# like a plain integer (e.g. "1"), by building it with a formula and then
# collapsing the formula to its computed (text) value via Paste Special
# Values. A straight ".Value = '1'" would be auto-coerced to the number 1,
# same as typing 1 into Excel directly - this keeps the User ID column as
# text like the rest of the report.
function Set-TextValue($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# --- Header row: the per-answer-option counters in D1:H1 shift right by one slot ---
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 1
$ws.Range("F1").Value = 2
$ws.Range("G1").Value = 3
$ws.Range("H1").Value = 4

# --- Row 2: User ID 1 (was 2), new session/orig IDs, answer (D2) unchanged ---
Set-TextValue $ws.Range("A2") "1"
$ws.Range("B2").Value = "0ee7e2E-3-83"
$ws.Range("C2").Value = "0ee7e2E-3-83"

# --- Row 3: User ID 2 (was 3), new session/orig IDs, answer changes to PROBABLY_NOT ---
Set-TextValue $ws.Range("A3") "2"
$ws.Range("B3").Value = "1iG-2I5c1-4-9"
$ws.Range("C3").Value = "1iG-2I5c1-4-9"
$ws.Range("E3").Value = "PROBABLY_NOT"

# --- Row 4: User ID 3 (was 4), new session/orig IDs, answer (F4) unchanged ("NO") ---
Set-TextValue $ws.Range("A4") "3"
$ws.Range("B4").Value = "2EI-4i2G-5-2-7"
$ws.Range("C4").Value = "2EI-4i2G-5-2-7"

# --- Row 5: User ID 4 (was 5), new session/orig IDs, answer changes to "NO" ---
Set-TextValue $ws.Range("A5") "4"
$ws.Range("B5").Value = "3gI-3a8g05-8"
$ws.Range("C5").Value = "3gI-3a8g05-8"
$ws.Range("G5").Value = "NO"

# --- Row 6: User ID 5 (was 1), now has its own session/orig IDs, answer becomes I_CANT_TELL ---
Set-TextValue $ws.Range("A6") "5"
$ws.Range("B6").Value = "4Ei0e-6g5-3-1"
$ws.Range("C6").Value = "4Ei0e-6g5-3-1"
$ws.Range("H6").Value = "I_CANT_TELL"

# --- Row 7: "Unanswered" summary counts change ---
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 6
$ws.Range("K7").Value = 6
$ws.Range("L7").Value = 6
$ws.Range("M7").Value = 6

# --- Column widths: re-fit the best-fit columns whose content width changed ---
$ws.Columns.Item(2).ColumnWidth = 12.918
$ws.Columns.Item(3).ColumnWidth = 12.918
$ws.Columns.Item(5).ColumnWidth = 14.418
$ws.Columns.Item(7).ColumnWidth = 2.7515
$ws.Columns.Item(8).ColumnWidth = 11.585
